# Modifications to model run: add an "infiltration_cmhr" column (W) to the
# units_m_day soils table, narrow column V, and update the active view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("units_m_day")
$ws.Activate()

# --- New column W: infiltration_cmhr -----------------------------------
# Copy the formatting of the last existing header cell (V1) onto the new
# header cell so the new column matches the surrounding header style,
# then write the header text and the data values for rows 2-11.
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("W1").Value = "infiltration_cmhr"
$ws.Range("W2:W11").Value = 1

# --- Column width for column V ------------------------------------------
$ws.Columns.Item(22).ColumnWidth = 4.6

# --- View / selection state ----------------------------------------------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B1").Select()
$win.FreezePanes = $true
$ws.Range("P22").Select()
